$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tijdsbestedingen na maandag: week 09-15/11/2015 changed from 4 to 7 hours
$ws.Range("B6").Formula = "=7"

$excel.Calculate()
